$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current layout (before edit):
#   row 179: 09-09-2024 / 943.18   (last data row)
#   row 180: <empty>
#   row 181: "Pie de página: Reporte generado automáticamente."
#   row 182: "Última actualización: 2024-09-08T07:38:50Z"
#
# Target layout (after edit):
#   row 179: 09-09-2024 / 943.18
#   row 180: 09-10-2024 / 946.22   (new)
#   row 181: 09-11-2024 / 948.85   (new)
#   row 182: <empty>
#   row 183: "Pie de página: Reporte generado automáticamente."
#   row 184: "Última actualización: 2024-09-10T21:15:11Z"

# Clear the old footer rows (181 and 182) so nothing stale lingers.
$ws.Range("A181").Value = ""
$ws.Range("A182").Value = ""

# Write the two new data rows. Use a leading apostrophe so Excel stores the
# date-like / numeric-like text as a literal string instead of auto-
# converting it to a date serial / number, then reset the style so no
# quote-prefix formatting sticks to the cell.
$ws.Range("A180").Value = "'09-10-2024"
$ws.Range("B180").Value = "'946.22"
$ws.Range("A181").Value = "'09-11-2024"
$ws.Range("B181").Value = "'948.85"
$ws.Range("A180:B181").Style = "Normal"

# Re-write the footer rows at their new location.
$ws.Range("A183").Value = "Pie de página: Reporte generado automáticamente."
$ws.Range("A184").Value = "Última actualización: 2024-09-10T21:15:11Z"
